$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: force text storage for numeric-looking price strings by
# temporarily marking the cell as Text before assignment, then restoring the
# default "Normal" style so no stray number-format style is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.367.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.687.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5454"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.06%  "
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06455"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07682"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.697.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.533"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5811"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008368"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.435.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.935"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.010"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.248"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.012"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.06%  "
$ws.Range("E25").Value = "  +5.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.862"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06316"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.65%  "
$ws.Range("E29").Value = "  +4.78%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.590"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.689"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6167"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.411"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.716"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.281"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.110.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8806"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.838.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000112"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.207"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.011"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05269"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4307"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.037"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.17%  "
